$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.945.71'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.905.93'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.30%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9981'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.8358'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +9.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.69'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.60%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9984'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3204'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.50%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '26.64'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.97%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07001'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08008'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.63%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7528'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.96%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.910.51'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.203'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.68'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.92%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.950.61'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.11'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.881'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '245.26'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007748'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.65%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.155.16'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.72%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9982'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9976'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.968'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1641'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +27.51%  '
$ws.Range("E26").Value = '  +1.60%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.224'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.15%  '
$ws.Range("E28").Value = '  +1.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.082'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.91%  '
$ws.Range("E30").Value = '  -1.88%  '
$ws.Range("E31").Value = '  -0.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.299'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.86%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05589'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.082'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.52%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.273'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.27%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7345'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.97%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.705'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01921'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.15%  '
$ws.Range("E39").Value = '  +0.32%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4427'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '72.27'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.980'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.85%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9969'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8375'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.893'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.29%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.587'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.67%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '100.82'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.90%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.735'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.88%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '988.07'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +9.63%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.061.13'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.94%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '36.23'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.45%  '
